$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 so the new header cells I1/J1 reuse the
# same cell style (s="1") as the other header cells instead of Excel
# allocating a brand-new style entry.
$ws.Range("H1:H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-9 for new columns I (I0) and J (IF)
$i0 = @(8, 9, 7, 9, 4, 7, 7, 5)
$if = @(9, 9, 9, 9, 7, 9, 8, 5)

for ($idx = 0; $idx -lt 8; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0[$idx]
    $ws.Cells.Item($row, 10).Value = $if[$idx]
}
